# Sample dataset update: demonstrates all cleaning operations.
# A new "PID" variable row is inserted after "ROWID" (renamed from "ID"),
# which shifts all subsequent variable rows down by one. Several numeric
# ranges / check results are refreshed, and two new variable rows
# ("Date of Treatment" content shifted down, plus a brand new
# "Total Care Cost" row) are appended at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Str($ref, $text) {
    $ws.Range($ref).Value = $text
}

function Set-Num($ref, $num) {
    $ws.Range($ref).Value = $num
}

function Clear-Cell($ref) {
    $ws.Range($ref).Value = ""
}

# ---- Row 2: ID -> ROWID ----
Set-Str "A2" "ROWID"
Set-Str "G2" "[1,29]"

# ---- Row 3: Name -> PID (new variable inserted) ----
Set-Str "A3" "PID"
Set-Str "B3" "Int64"
Set-Str "C3" "numeric"
Set-Str "G3" "[1,20]"

# ---- Row 4: Gender -> Name ----
Set-Str "A4" "Name"
Clear-Cell "I4"

# ---- Row 5: Age -> Gender ----
Set-Str "A5" "Gender"
Set-Str "B5" "string"
Set-Str "C5" "string"
Clear-Cell "G5"
Set-Str "I5" "male, female"

# ---- Row 6: AgeMonths -> Age ----
Set-Str "A6" "Age"
Set-Str "G6" "[16,39]"

# ---- Row 7: Height -> AgeMonths ----
Set-Str "A7" "AgeMonths"
Set-Str "B7" "Int64"
Set-Str "G7" "[194,476]"

# ---- Row 8: Weight -> Height ----
Set-Str "A8" "Height"
Set-Num "E8" 2
Set-Num "F8" 0.1
Set-Str "G8" "[-96.2,190.1]"
Set-Str "I8" "[0,200]"
Set-Str "J8" "check_true"
$ws.Range("L8").NumberFormat = "@"
Set-Str "L8" "18"

# ---- Row 9: BMI -> Weight ----
Set-Str "A9" "Weight"
Set-Num "E9" 2
Set-Num "F9" 0.1
Set-Str "G9" "[-49.9,88.7]"
Set-Str "I9" "[0,150]"
Set-Str "J9" "check_true"
$ws.Range("L9").NumberFormat = "@"
Set-Str "L9" "19"

# ---- Row 10: BMICatUnder20yrs -> BMI ----
Set-Str "A10" "BMI"
Set-Str "B10" "Float64"
Set-Str "C10" "numeric"
Set-Num "E10" 3
Set-Num "F10" 0.15
Set-Str "G10" "[-19.89,95.85]"
Clear-Cell "I10"
Set-Str "K10" "check_false"
Set-Str "M10" "N.A."

# ---- Row 11: BMI_WHO -> BMICatUnder20yrs ----
Set-Str "A11" "BMICatUnder20yrs"
Set-Str "I11" "UnderWeight; NormWeight; OverWeight; Obese"
Set-Str "M11" "1,3,6,7,8,10,11,12,15,16,17,18,19"

# ---- Row 12: Date of Birth -> BMI_WHO ----
Set-Str "A12" "BMI_WHO"
Set-Str "B12" "string"
Set-Str "C12" "string"
Set-Str "I12" "12.0_18.5; 18.5_to_24.9; 25.0_to_29.9; 30.0_plus"
Set-Str "K12" "check_true"
Set-Str "M12" "12,16,17"

# ---- Row 13: Date of First Visit -> Date of Birth ----
Set-Str "A13" "Date of Birth"

# ---- Row 14: Date of Diagnosis -> Date of First Visit ----
Set-Str "A14" "Date of First Visit"

# ---- Row 15: Date of Treatment -> Date of Diagnosis ----
Set-Str "A15" "Date of Diagnosis"
Set-Num "E15" 0
Set-Num "F15" 0

# ---- Row 16 (new): Date of Treatment ----
Set-Str "A16" "Date of Treatment"
Set-Str "B16" "object"
Set-Str "C16" "date"
Set-Str "D16" "Matched"
Set-Num "E16" 10
Set-Num "F16" 0.5
Set-Str "H16" "N.A."
Set-Str "I16" "ddd, dd mmmm yy"
Set-Str "J16" "check_false"
Set-Str "K16" "check_false"
Set-Str "L16" "N.A."
Set-Str "M16" "N.A."

# ---- Row 17 (new): Total Care Cost ----
Set-Str "A17" "Total Care Cost"
Set-Str "B17" "string"
Set-Str "C17" "string"
Set-Str "D17" "Matched"
Set-Num "E17" 0
Set-Num "F17" 0
Set-Str "H17" "N.A."
Set-Str "J17" "check_false"
Set-Str "K17" "check_false"
Set-Str "L17" "N.A."
Set-Str "M17" "N.A."

# Column A on the variable-name column uses a bold/centered/bordered style
# (same as the other rows in that column); copy it onto the two new rows.
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16:A17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
